$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (e.g. "584.62", "0.100") that must stay
# literal text like the source data, not be auto-converted to numbers by Excel.
# Forcing each target cell to Text format ("@") right before the assignment
# keeps the written value a string.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.093.45"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.267.97"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.62"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "184.24"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.600"
$ws.Range("E8").Value = "  +0.38%  "
$ws.Range("E9").Value = "  -1.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.63"
$ws.Range("E10").Value = "  -0.44%  "
$ws.Range("E11").Value = "  -2.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.832.82"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("E13").Value = "  +1.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "68.080.49"
$ws.Range("E14").Value = "  -0.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.35"
$ws.Range("E15").Value = "  -3.01%  "
$ws.Range("E16").Value = "  -1.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.256.97"
$ws.Range("E17").Value = "  -0.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.72"
$ws.Range("E18").Value = "  -2.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.26"
$ws.Range("E19").Value = "  -2.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "416.73"
$ws.Range("E20").Value = "  +5.97%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.53"
$ws.Range("E21").Value = "  -1.99%  "
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.04"
$ws.Range("E23").Value = "  -0.29%  "
$ws.Range("E24").Value = "  -2.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000117"
$ws.Range("E25").Value = "  -2.19%  "
$ws.Range("E26").Value = "  -0.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.34"
$ws.Range("E27").Value = "  -4.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.997"
$ws.Range("E28").Value = "  -0.36%  "
$ws.Range("E29").Value = "  -1.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.63"
$ws.Range("E30").Value = "  -1.16%  "
$ws.Range("E31").Value = "  -5.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.84"
$ws.Range("E32").Value = "  -4.53%  "
$ws.Range("E33").Value = "  -3.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "164.42"
$ws.Range("E34").Value = "  +0.89%  "
$ws.Range("E35").Value = "  -4.83%  "
$ws.Range("E36").Value = "  -3.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.66"
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("E38").Value = "  -3.86%  "
$ws.Range("E39").Value = "  -3.41%  "
$ws.Range("E40").Value = "  -4.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.627.93"
$ws.Range("E41").Value = "  -0.94%  "
$ws.Range("E42").Value = "  -2.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.40"
$ws.Range("E43").Value = "  -3.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "335.60"
$ws.Range("E44").Value = "  -0.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "24.16"
$ws.Range("E45").Value = "  -4.96%  "
$ws.Range("E46").Value = "  -3.09%  "
$ws.Range("E47").Value = "  -2.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.983"
$ws.Range("E48").Value = "  -0.60%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.100"
$ws.Range("E49").Value = "  -1.65%  "
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "30.58"
$ws.Range("E51").Value = "  -3.07%  "
